# The source data for rows 57-60 ("Artfynd" sheet) got re-synced upstream:
# the four observation records now appear in a different order/rotation.
# Row 57 <- old row 58's data, Row 58 <- old row 60's data,
# Row 59 <- old row 57's data, Row 60 <- old row 59's data.
# Below we just write out each row's new field values explicitly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $val) {
    # Force the value to be stored as text even when it looks numeric
    # (mirrors the source column being formatted/typed as Text in Excel).
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

# Row 57
$ws.Cells.Item(57, 1).Value = 112145544
$ws.Cells.Item(57, 2).Value = 96735
$ws.Cells.Item(57, 4).Value = "VU"
$ws.Cells.Item(57, 5).Value = 220787
$ws.Cells.Item(57, 6).Value = "Knärot"
$ws.Cells.Item(57, 7).Value = "Goodyera repens"
$ws.Cells.Item(57, 8).Value = "(L.) R. Br."
Set-TextCell $ws.Cells.Item(57, 9) "17"
$ws.Cells.Item(57, 10).Value = "plantor/tuvor"
$ws.Cells.Item(57, 13).ClearContents()
$ws.Cells.Item(57, 16).Value = "Bennarby, Upl"
$ws.Cells.Item(57, 17).Value = 653024
$ws.Cells.Item(57, 18).Value = 6675364
$ws.Cells.Item(57, 19).Value = 4
$ws.Cells.Item(57, 26).Value = "11:03"
$ws.Cells.Item(57, 28).Value = "11:04"

# Row 58
$ws.Cells.Item(58, 1).Value = 112145539
$ws.Cells.Item(58, 2).Value = 90166
$ws.Cells.Item(58, 4).Value = "LC"
$ws.Cells.Item(58, 5).Value = 1339
$ws.Cells.Item(58, 6).Value = "Brandticka"
$ws.Cells.Item(58, 7).Value = "Pycnoporellus fulgens"
$ws.Cells.Item(58, 8).Value = "(Fr.) Donk"
Set-TextCell $ws.Cells.Item(58, 9) ""
$ws.Cells.Item(58, 10).Value = ""
$ws.Cells.Item(58, 12).ClearContents()
$ws.Cells.Item(58, 17).Value = 652997
$ws.Cells.Item(58, 18).Value = 6675310
$ws.Cells.Item(58, 26).Value = "10:44"
$ws.Cells.Item(58, 28).Value = "10:44"

# Row 59
$ws.Cells.Item(59, 1).Value = 112145535
$ws.Cells.Item(59, 2).Value = 56575
$ws.Cells.Item(59, 4).Value = "NT"
$ws.Cells.Item(59, 5).Value = 103021
$ws.Cells.Item(59, 6).Value = "Talltita"
$ws.Cells.Item(59, 7).Value = "Poecile montanus"
$ws.Cells.Item(59, 8).Value = "(Conrad von Baldenstein, 1827)"
Set-TextCell $ws.Cells.Item(59, 9) "1"
$ws.Cells.Item(59, 10).ClearContents()
$ws.Cells.Item(59, 13).Value = "permanent revir"
$ws.Cells.Item(59, 16).Value = "Smigruvan, Upl"
$ws.Cells.Item(59, 17).Value = 653012
$ws.Cells.Item(59, 18).Value = 6675152
$ws.Cells.Item(59, 19).Value = 84
$ws.Cells.Item(59, 26).Value = "10:10"
$ws.Cells.Item(59, 28).Value = "10:10"

# Row 60
$ws.Cells.Item(60, 1).Value = 112145545
$ws.Cells.Item(60, 2).Value = 96735
$ws.Cells.Item(60, 4).Value = "VU"
$ws.Cells.Item(60, 5).Value = 220787
$ws.Cells.Item(60, 6).Value = "Knärot"
$ws.Cells.Item(60, 7).Value = "Goodyera repens"
$ws.Cells.Item(60, 8).Value = "(L.) R. Br."
Set-TextCell $ws.Cells.Item(60, 9) "14"
$ws.Cells.Item(60, 10).Value = "plantor/tuvor"
$ws.Cells.Item(60, 12).Value = ""
$ws.Cells.Item(60, 17).Value = 653038
$ws.Cells.Item(60, 18).Value = 6675341
$ws.Cells.Item(60, 26).Value = "11:07"
$ws.Cells.Item(60, 28).Value = "11:08"
